# Add a new row to the "Completed" reading list for
# "The Rise and Fall of the Dinosaurs"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$row = 49

$ws.Cells.Item($row, 1).Value = "The Rise and Fall of the Dinosaurs"
$ws.Cells.Item($row, 2).Value = "Steve Brusatte;Patrick Lawlor"
$ws.Cells.Item($row, 3).Value = 43921
$ws.Cells.Item($row, 4).Value = 43922
$ws.Cells.Item($row, 5).Value = "dinosaurs;science;history"
$ws.Cells.Item($row, 6).Value = "Audio"
$ws.Cells.Item($row, 7).Value = "10 Hours 12 Mins"

# Reuse the same date style already used in the row above (column C/D)
# instead of creating a brand new number format.
$ws.Range("C48:D48").Copy()
$ws.Range("C49:D49").PasteSpecial(-4122)

$ws.Activate()
$ws.Range("A50").Select()
